$d = $word.ActiveDocument

function Replace-All($findText, $replaceText) {
    $rng = $d.Content
    $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
}

function Replace-Once($findText, $replaceText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $replaceText
    }
}

# Simple "placeholder + trailing room number" cells: the whole run is the
# underscore/comma placeholder followed by a room number, e.g.
# "________________,  213a/4" -> "Cineva  213a/4" (4 occurrences)
Replace-All "________________,  213a/4" "Cineva  213a/4"

# "______________,  350/4" -> "Cineva  350/4" (4 occurrences total: sz24
# yellow + sz32 yellow cells)
Replace-All "______________,  350/4" "Cineva  350/4"

# "_________________, 350/4" -> "Cineva 350/4" (4 occurrences, lightGray cells)
Replace-All "_________________, 350/4" "Cineva 350/4"

# Cells where the room number was already split into its own run
# ("113"/"214"/"423") -- the placeholder text still reads the same when
# traversed as plain text across run boundaries, so Find still locates it.
Replace-Once "________________,  113/4" "Cineva  113/4"
Replace-Once "________________,  214/4" "Cineva  214/4"

# "_________________, 423/4" occurs twice, unchanged apart from the
# placeholder -> Cineva swap.
Replace-All "_________________, 423/4" "Cineva 423/4"
